$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Phone" field (row 4)
$ws.Range("A4").Value = "Phone"
$ws.Range("B4").Value = 712345678

# Add the new "Description" field (row 5)
$ws.Range("A5").Value = "Description"
$ws.Range("B5").Value = "Test RPA"

# Update selection to match the post-edit state
$ws.Range("C5").Select()
